$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A:C width change (target stored width 37.85546875 chars;
# COM ColumnWidth snaps to a pixel grid, so 37.0 is the input that
# round-trips to the closest achievable stored width, 37.833333...)
$ws.Range("A1:C1").EntireColumn.ColumnWidth = 37

# Row height changes
$ws.Rows.Item(4).RowHeight = 28.5
$ws.Rows.Item(5).RowHeight = 27.75
$ws.Rows.Item(6).RowHeight = 28.5
$ws.Rows.Item(7).RowHeight = 27
$ws.Rows.Item(8).RowHeight = 18.75
$ws.Rows.Item(10).RowHeight = 18.75
$ws.Rows.Item(11).RowHeight = 17.25

# Value changes
$ws.Range("Q4").Value = 109
$ws.Range("Q7").Value = 12685.1
$ws.Range("P8").Value = 478225.6
$ws.Range("Q8").Value = 559503.6
$ws.Range("Q9").Value = 131.9
$ws.Range("Q10").Value = 3384.8
$ws.Range("Q11").Value = 12517.9
